$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D values are written as text, preserving exact formatting
# (e.g. trailing zeros, thousands-dot grouping) matching the original inlineStr cells.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.747.23"
$ws.Range("E2").Value = "  -2.43%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.564.28"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "206.14"
$ws.Range("E5").Value = "  -0.89%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.488"
$ws.Range("E6").Value = "  -2.27%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "21.87"
$ws.Range("E8").Value = "  -0.46%  "
$ws.Range("E9").Value = "  -0.40%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0583"
$ws.Range("E10").Value = "  -1.22%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0861"
$ws.Range("E11").Value = "  -0.41%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.788.67"
$ws.Range("E12").Value = "  +0.23%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.570.63"
$ws.Range("E13").Value = "  +0.58%  "
$ws.Range("E14").Value = "  -2.12%  "
$ws.Range("E15").Value = "  -0.34%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "26.802.48"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.49"
$ws.Range("E17").Value = "  -3.05%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "213.95"
$ws.Range("E18").Value = "  +0.79%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.33"
$ws.Range("E19").Value = "  +1.51%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0676"
$ws.Range("E20").Value = "  -1.57%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.08"
$ws.Range("E22").Value = "  -0.14%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.31"
$ws.Range("E23").Value = "  -1.72%  "
$ws.Range("E24").Value = "  -0.86%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "152.55"
$ws.Range("E25").Value = "  -0.04%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.73"
$ws.Range("E26").Value = "  +0.94%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "14.90"
$ws.Range("E27").Value = "  -0.21%  "
$ws.Range("E28").Value = "  +0.08%  "
$ws.Range("E29").Value = "  -1.32%  "
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.11"
$ws.Range("E30").Value = "  -3.32%  "
$ws.Range("B31").Value = "Hedera"
$ws.Range("C31").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0462"
$ws.Range("E31").Value = "  -1.38%  "
$ws.Range("E32").Value = "  -1.41%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.384.88"
$ws.Range("E33").Value = "  +1.29%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.92"
$ws.Range("E34").Value = "  -0.94%  "
$ws.Range("E35").Value = "  +0.97%  "
$ws.Range("E36").Value = "  -0.82%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.925"
$ws.Range("E37").Value = "  -3.08%  "
$ws.Range("E38").Value = "  -2.37%  "
$ws.Range("E39").Value = "  -0.92%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.815"
$ws.Range("E40").Value = "  -0.27%  "
$ws.Range("E41").Value = "  +0.05%  "
$ws.Range("E42").Value = "  +1.53%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.34"
$ws.Range("E43").Value = "  +1.81%  "
$ws.Range("E44").Value = "  -0.59%  "
$ws.Range("E45").Value = "  +0.92%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "63.21"
$ws.Range("E46").Value = "  -0.72%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.701.14"
$ws.Range("E47").Value = "  +0.22%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "85.37"
$ws.Range("E48").Value = "  +0.11%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0₇0986"
$ws.Range("E49").Value = "  -0.12%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0950"
$ws.Range("E50").Value = "  -0.32%  "
$ws.Range("E51").Value = "  -0.75%  "
